# Update cryptocurrency price/volume data on sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Straightforward text-cell updates ---
# (these values survive Excel's auto-number-detection as text unchanged)
$ws.Range('D2').Value = '29.392.88'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '1.848.07'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('D5').Value = '240.28'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').Value = '0.6303'
$ws.Range('D8').Value = '0.07607'
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('D9').Value = '0.2937'
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('D10').Value = '24.57'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('D11').Value = '0.07742'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = '1.882.78'
$ws.Range('E12').Value = '  -5.15%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E13').Value = '  +9.99%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '5.008'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = '83.66'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '2.153.54'
$ws.Range('E17').Value = '  -4.89%  '
$ws.Range('D18').Value = '6.195'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = '29.434.51'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '229.42'
$ws.Range('E20').Value = '  -0.99%  '
$ws.Range('D21').Value = '12.45'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '7.527'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('D25').Value = '157.23'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('D26').Value = '0.1402'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').Value = '8.364'
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('D28').Value = '17.65'
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('D29').Value = '1.463'
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('D30').Value = '1.299'
$ws.Range('E30').Value = '  +2.99%  '
$ws.Range('D31').Value = '0.05592'
$ws.Range('E31').Value = '  -2.02%  '
$ws.Range('D32').Value = '4.115'
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('D33').Value = '4.038'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '1.859'
$ws.Range('E34').Value = '  +0.39%  '
$ws.Range('D35').Value = '1.158'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '0.7128'
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('E37').Value = '  -0.43%  '
$ws.Range('D38').Value = '1.246.33'
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').Value = '2.778'
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('E41').Value = '  +5.54%  '
$ws.Range('D42').Value = '0.9037'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '1.001'
$ws.Range('D44').Value = '101.84'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').Value = '65.99'
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('E46').Value = '  +1.46%  '
$ws.Range('D47').Value = '7.141'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('D50').Value = '8.974'
$ws.Range('E50').Value = '  -1.95%  '
$ws.Range('D51').Value = '0.1122'
$ws.Range('E51').Value = '  -0.47%  '

# --- Updates for values that Excel would otherwise coerce to a Number ---
# (e.g. "0.00001096", "0.6800", "6.430", "1.690" all round-trip cleanly through
#  General number formatting, so a plain .Value/.Formula assignment would silently
#  turn them into numeric cells and lose the original text formatting/precision).
# Work around this by building the literal text in a scratch cell via a formula
# (text formulas are never re-interpreted as numbers), then copy/paste-special
# VALUES ONLY into the destination cell, which keeps it a genuine text cell.
$scratchRow = 100
$scratch = $ws.Cells.Item($scratchRow, 1)
$scratch.NumberFormat = "@"

$scratch.Formula = '="0.00001096"'
$scratch.Copy()
$ws.Range('D13').PasteSpecial(-4163)

$scratch.Formula = '="0.6800"'
$scratch.Copy()
$ws.Range('D15').PasteSpecial(-4163)

$scratch.Formula = '="6.430"'
$scratch.Copy()
$ws.Range('D41').PasteSpecial(-4163)

$scratch.Formula = '="1.690"'
$scratch.Copy()
$ws.Range('D49').PasteSpecial(-4163)

# Remove the scratch row entirely so it does not affect the used range / dimension
$scratch.EntireRow.Delete()

$excel.CutCopyMode = 0

